$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated with new TPM-derived receptor expression values and
# all downstream specificity / edge-weight figures recomputed from them.
$ws.Range("M2").Value = 0.6746743333333334
$ws.Range("N2").Value = 2.024023
$ws.Range("O2").Value = 0.07069047851636343
$ws.Range("P2").Value = 0.07069047851636343
$ws.Range("Q2").Value = 0.6953755907944446
$ws.Range("R2").Value = 6.25838031715
$ws.Range("S2").Value = 0.07069047851636343
$ws.Range("T2").Value = 0.07069047851636343

# Row 3 - only the specificity / edge-weight columns shift, because they
# are normalised against row 2's new receptor expression values.
$ws.Range("O3").Value = 0.4692497642600617
$ws.Range("P3").Value = 0.4692497642600616
$ws.Range("Q3").Value = 4.615965811816667
$ws.Range("R3").Value = 41.54369230635
$ws.Range("S3").Value = 0.4692497642600617
$ws.Range("T3").Value = 0.4692497642600616

# Row 4 - same normalisation effect as row 3.
$ws.Range("O4").Value = 0.460059757223575
$ws.Range("P4").Value = 0.460059757223575
$ws.Range("Q4").Value = 4.52556457665
$ws.Range("R4").Value = 40.73008118985
$ws.Range("S4").Value = 0.460059757223575
$ws.Range("T4").Value = 0.460059757223575
